$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Row, $Col, $Text) {
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

# Simple price/volume updates (rows 2-38, 41-43, 45-47, 49-50)
Set-TextCell 2 4 '30.925.32'
$ws.Cells.Item(2, 5).Value = '  +0.64%  '
Set-TextCell 3 4 '1.918.06'
$ws.Cells.Item(3, 5).Value = '  +1.18%  '
Set-TextCell 4 4 '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
Set-TextCell 5 4 '239.60'
$ws.Cells.Item(5, 5).Value = '  -3.45%  '
Set-TextCell 6 4 '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.11%  '
Set-TextCell 7 4 '0.4917'
$ws.Cells.Item(7, 5).Value = '  -0.48%  '
Set-TextCell 8 4 '0.2972'
$ws.Cells.Item(8, 5).Value = '  +0.04%  '
Set-TextCell 9 4 '0.06781'
$ws.Cells.Item(9, 5).Value = '  -0.59%  '
Set-TextCell 10 4 '1.911.12'
$ws.Cells.Item(10, 5).Value = '  +0.81%  '
Set-TextCell 11 4 '17.07'
$ws.Cells.Item(11, 5).Value = '  -1.28%  '
Set-TextCell 12 4 '0.07309'
$ws.Cells.Item(12, 5).Value = '  +0.59%  '
Set-TextCell 13 4 '5.152'
$ws.Cells.Item(13, 5).Value = '  +0.65%  '
Set-TextCell 14 4 '90.23'
$ws.Cells.Item(14, 5).Value = '  -2.36%  '
Set-TextCell 15 4 '0.6750'
$ws.Cells.Item(15, 5).Value = '  -0.82%  '
Set-TextCell 16 4 '30.895.93'
$ws.Cells.Item(16, 5).Value = '  +0.61%  '
Set-TextCell 17 4 '0.000007958'
$ws.Cells.Item(17, 5).Value = '  -0.40%  '
Set-TextCell 18 4 '13.46'
$ws.Cells.Item(18, 5).Value = '  +0.99%  '
Set-TextCell 19 4 '1.001'
$ws.Cells.Item(19, 5).Value = '  +0.03%  '
Set-TextCell 20 4 '2.167.15'
$ws.Cells.Item(20, 5).Value = '  +1.30%  '
Set-TextCell 21 4 '1.001'
$ws.Cells.Item(21, 5).Value = '  +0.05%  '
Set-TextCell 22 4 '5.180'
$ws.Cells.Item(22, 5).Value = '  +6.54%  '
Set-TextCell 23 4 '207.99'
$ws.Cells.Item(23, 5).Value = '  +7.50%  '
Set-TextCell 24 4 '6.247'
$ws.Cells.Item(24, 5).Value = '  +2.65%  '
Set-TextCell 25 4 '9.683'
$ws.Cells.Item(25, 5).Value = '  +2.54%  '
Set-TextCell 26 4 '158.27'
$ws.Cells.Item(26, 5).Value = '  +1.52%  '
Set-TextCell 27 4 '18.93'
$ws.Cells.Item(27, 5).Value = '  -1.77%  '
Set-TextCell 28 4 '1.980'
$ws.Cells.Item(28, 5).Value = '  +2.87%  '
Set-TextCell 29 4 '1.422'
$ws.Cells.Item(29, 5).Value = '  +1.30%  '
Set-TextCell 30 4 '4.329'
$ws.Cells.Item(30, 5).Value = '  -0.79%  '
Set-TextCell 31 4 '0.09192'
$ws.Cells.Item(31, 5).Value = '  +1.94%  '
Set-TextCell 32 4 '4.072'
$ws.Cells.Item(32, 5).Value = '  +0.76%  '
Set-TextCell 33 4 '0.05174'
$ws.Cells.Item(33, 5).Value = '  -0.77%  '
Set-TextCell 34 4 '0.7542'
$ws.Cells.Item(34, 5).Value = '  +0.94%  '
Set-TextCell 35 4 '1.127'
$ws.Cells.Item(35, 5).Value = '  -0.06%  '
Set-TextCell 36 4 '2.718'
$ws.Cells.Item(36, 5).Value = '  -0.76%  '
Set-TextCell 37 4 '0.01860'
$ws.Cells.Item(37, 5).Value = '  -0.27%  '
Set-TextCell 38 4 '2.736'
$ws.Cells.Item(38, 5).Value = '  +1.99%  '
Set-TextCell 41 4 '0.4522'
$ws.Cells.Item(41, 5).Value = '  +1.75%  '
Set-TextCell 42 4 '107.81'
$ws.Cells.Item(42, 5).Value = '  +1.24%  '
Set-TextCell 43 4 '5.908'
$ws.Cells.Item(43, 5).Value = '  +2.19%  '
Set-TextCell 45 4 '0.1407'
$ws.Cells.Item(45, 5).Value = '  +4.39%  '
Set-TextCell 46 4 '7.741'
$ws.Cells.Item(46, 5).Value = '  +0.61%  '
Set-TextCell 47 4 '66.41'
$ws.Cells.Item(47, 5).Value = '  +14.03%  '
Set-TextCell 49 4 '0.4110'
$ws.Cells.Item(49, 5).Value = '  +3.70%  '
Set-TextCell 50 4 '0.05951'
$ws.Cells.Item(50, 5).Value = '  +1.46%  '

# Row 44: only Volume(1h) changes
$ws.Cells.Item(44, 5).Value = '  +0.95%  '

# Rows with swapped coin identity (name/link) plus price/volume updates
$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 39 4 '2.110'
$ws.Cells.Item(39, 5).Value = '  -2.71%  '

$ws.Cells.Item(40, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 40 4 '0.9260'
$ws.Cells.Item(40, 5).Value = '  -1.97%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 48 4 '9.029'
$ws.Cells.Item(48, 5).Value = '  +3.60%  '

$ws.Cells.Item(51, 2).Value = 'Elrond'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 51 4 '35.00'
$ws.Cells.Item(51, 5).Value = '  +4.01%  '
